# Apply updated cryptocurrency price/volume data to the worksheet.
# Numeric-looking price strings (e.g. "221.31") must be forced to remain
# plain text (matching the original inlineStr cells) instead of being
# auto-coerced to numbers by COM. We do this by temporarily switching the
# cell to a text NumberFormat, assigning the value, then restoring the
# cell's original Style so no visible formatting changes remain.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.127.10'
$ws.Range("E2").Value = '  -1.40%  '
$ws.Range("D3").Value = '1.782.52'
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("E4").Value = '  +0.23%  '
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '221.31'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  -1.86%  '
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.548'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  -1.17%  '
$ws.Range("E7").Value = '  +0.29%  '
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.46'
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = '  -4.89%  '
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.287'
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = '  +0.77%  '
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0708'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  +6.40%  '
$ws.Range("E11").Value = '  -0.89%  '
$ws.Range("D12").Value = '2.052.23'
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").Value = '1.796.61'
$ws.Range("E13").Value = '  -0.26%  '
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.61'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  -4.43%  '
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.624'
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = '  -2.86%  '
$ws.Range("D16").Value = '34.024.53'
$ws.Range("E16").Value = '  -1.44%  '
$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.21'
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = '  -1.74%  '
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.89'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  -2.56%  '
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.14'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  -4.62%  '
$ws.Range("D20").Value = '0.0₃0775'
$ws.Range("E20").Value = '  +3.02%  '
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  +0.34%  '
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.62'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  +1.31%  '
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.08'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  -3.77%  '
$ws.Range("E24").Value = '  -0.15%  '
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.71'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  +0.03%  '
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.24'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  -1.78%  '
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.02'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  -1.31%  '
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.112'
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = '  -2.16%  '
$ws.Range("E29").Value = '  +0.60%  '
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0518'
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = '  +0.22%  '
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.67'
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = '  -3.24%  '
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.19'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  +0.07%  '
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.50'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  -2.97%  '
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.83'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  -4.37%  '
$ws.Range("D35").Value = '1.396.47'
$ws.Range("E35").Value = '  -4.39%  '
$ws.Range("E36").Value = '  -0.84%  '
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.626'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  -1.48%  '
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0185'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  -2.64%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.76'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  -3.44%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.928'
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = '  +3.15%  '
$ws.Range("B41").Value = 'HuobiToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.35'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  +1.75%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '79.22'
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = '  -4.86%  '
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.11'
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = '  +0.67%  '
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0493'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  -2.95%  '
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.05'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  +0.57%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.83'
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = '  -1.28%  '
$ws.Range("D47").Value = '1.941.09'
$ws.Range("E47").Value = '  -0.87%  '
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '105.46'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  +5.16%  '
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.997'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  -0.09%  '
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.66'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  -3.75%  '
$ws.Range("D51").Value = '0.0₆0116'
$ws.Range("E51").Value = '  -0.17%  '
